# Add a new "type" column (F) classifying each variable as categorical
# ("cat") or continuous ("con"), matching the upstream commit that added
# this column to data/obtainWeights.xlsx.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header
$ws.Range("F1").Value = "type"

# Row -> type ("cat" = categorical, "con" = continuous)
$ws.Range("F2").Value  = "cat"   # sex
$ws.Range("F3").Value  = "con"   # age
$ws.Range("F4").Value  = "con"   # education_age
$ws.Range("F5").Value  = "cat"   # alcfrequency
$ws.Range("F6").Value  = "cat"   # smoking_status
$ws.Range("F7").Value  = "cat"   # income
$ws.Range("F8").Value  = "cat"   # household_size
$ws.Range("F9").Value  = "cat"   # employment_status
$ws.Range("F10").Value = "con"   # bmi
$ws.Range("F11").Value = "cat"   # bmi_cat
$ws.Range("F12").Value = "cat"   # overallhealth
$ws.Range("F13").Value = "con"   # height
$ws.Range("F14").Value = "cat"   # urbanisation
$ws.Range("F15").Value = "con"   # weight
$ws.Range("F16").Value = "cat"   # assessment_center
$ws.Range("F17").Value = "cat"   # ethnic_background
$ws.Range("F18").Value = "cat"   # education_degree

# A subset of the new column's cells carry a plain black font (as opposed
# to the workbook default "theme" font color) in the source file.
$blackFontRows = @(11, 12, 14, 15, 16, 17, 18)
foreach ($r in $blackFontRows) {
    $ws.Range("F$r").Font.Color = 0
}

# Keep the selection / used range in sync with the appended column, same
# as Excel would after the edit (cursor parked one row below the table).
[void]$ws.Range("F19").Select()
